# Commit: "Fruta / hortaliza, semanal"
# This workbook records daily/weekly produce price observations, one row per
# (date, variety) combination, ordered chronologically. This edit adds a new
# weekly batch of 2 observations (date 2022-02-11) for "Pimiento" at the top
# of the existing data block (rows 153-154), pushing all subsequent rows
# down by two positions. The sheet's used range grows from A1:R242 to
# A1:R244 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 153-154; everything currently at row 153 onward
# shifts down by two rows (Excel also copies the row-153 formatting, e.g.
# the date number format on column D, down onto the freshly inserted rows).
$ws.Rows("153:154").Insert()

# --- New row 153: Zafiro rojo, 2022-02-11 ---
$ws.Cells.Item(153, 1).Value = 11
$ws.Cells.Item(153, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(153, 3).Value = "Bíobío"
$ws.Cells.Item(153, 4).Value = "2022-02-11"
$ws.Cells.Item(153, 5).Value = 8
$ws.Cells.Item(153, 6).Value = 100112002
$ws.Cells.Item(153, 7).Value = "Pimiento"
$ws.Cells.Item(153, 8).Value = "Zafiro rojo"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 180
$ws.Cells.Item(153, 11).Value = 16000
$ws.Cells.Item(153, 12).Value = 18000
$ws.Cells.Item(153, 13).Value = 16889
$ws.Cells.Item(153, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(153, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(153, 16).Value = 1126
$ws.Cells.Item(153, 17).Value = 15
$ws.Cells.Item(153, 18).Value = "Hortaliza"

# --- New row 154: Zafiro verde, 2022-02-11 ---
$ws.Cells.Item(154, 1).Value = 11
$ws.Cells.Item(154, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(154, 3).Value = "Bíobío"
$ws.Cells.Item(154, 4).Value = "2022-02-11"
$ws.Cells.Item(154, 5).Value = 8
$ws.Cells.Item(154, 6).Value = 100112002
$ws.Cells.Item(154, 7).Value = "Pimiento"
$ws.Cells.Item(154, 8).Value = "Zafiro verde"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 150
$ws.Cells.Item(154, 11).Value = 13000
$ws.Cells.Item(154, 12).Value = 14000
$ws.Cells.Item(154, 13).Value = 13333
$ws.Cells.Item(154, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(154, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(154, 16).Value = 889
$ws.Cells.Item(154, 17).Value = 15
$ws.Cells.Item(154, 18).Value = "Hortaliza"
